# Cashflows.xlsx — sign-flip revisions + new overview row
#
# 1. On the "Cashflow" sheet, rows 2-5 (CL amount type) and rows 6-9 (PR
#    amount type) across columns F:BN have their numeric signs flipped
#    (cash inflows/outflows corrected).
# 2. The old helper/formatting row 17 is removed and a small new
#    (currently-empty) row 11 is started, shifting the sheet's used range
#    from A1:BV17 down to A1:BV11.
# 3. The current selection on "Cashflow" is left on the row below the new
#    data block (A12:XFD20), matching where the editor was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cashflow")

# --- 1. Flip the sign of every value in F2:BN9 -----------------------------
$rng = $ws.Range("F2:BN9")
foreach ($cell in $rng) {
    $v = $cell.Value()
    if ($v -ne $null) {
        $cell.Value = -$v
    }
}

# --- 2. Remove the old formatting-only row 17, start a fresh row 11 -------
$ws.Rows.Item(17).Delete()
$ws.Range("F11").NumberFormat = "0.00"

# --- 3. Update the active selection on the Cashflow sheet ------------------
$activated = $ws.Activate()
$selected = $ws.Range("A12:XFD20").Select()
